# Applies the changes described in the commit:
#  - "Anonymize fedcore": rename the approach label "fedcore" -> "approach"
#    on both comparison sheets.
#  - Add a thin top/bottom border to the first two cells of each merged
#    "original/approach/change" header group, and a thin top/bottom/right
#    border to the last cell of each group (closing the box around the
#    merged header cell).
#  - Drop the stray empty inline-string cell at G5 on the
#    computational_comparison sheet.

$wb = $excel.ActiveWorkbook

function Set-GroupBorders($ws, [string]$midAddr, [string]$rightAddr) {
    # Middle cell of the merged group (e.g. C1 / F1): thin top + bottom.
    $ws.Range($midAddr).ClearFormats()
    $ws.Range($midAddr).Borders(8).Weight = 2   # xlEdgeTop
    $ws.Range($midAddr).Borders(9).Weight = 2   # xlEdgeBottom

    # Last cell of the merged group (e.g. D1 / G1): thin top + bottom + right.
    $ws.Range($rightAddr).ClearFormats()
    $ws.Range($rightAddr).Borders(8).Weight = 2   # xlEdgeTop
    $ws.Range($rightAddr).Borders(9).Weight = 2   # xlEdgeBottom
    $ws.Range($rightAddr).Borders(10).Weight = 2  # xlEdgeRight
}

# --- Sheet 1: quality_comparison --------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-GroupBorders $ws1 "C1" "D1"

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-GroupBorders $ws2 "C1" "D1"
Set-GroupBorders $ws2 "F1" "G1"

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell.
$ws2.Range("G5").ClearContents()
